$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns (E,F,G) before the old "Faculty names" column (which shifts to H)
$ws.Range("E1:G1").EntireColumn.Insert()

# Header row values for the new columns
$ws.Range("E1").Value = "Wed"
$ws.Range("F1").Value = "Thu"
$ws.Range("G1").Value = "Fri"

# Bold the whole header row (uses the default 12pt font, just bold)
$ws.Range("A1:H1").Font.Bold = $true

# Fill in the day-availability columns with 1 (available) by default
$ws.Range("E2:G27").Value = 1

# Row 6 (Yar, Tasha) is only available Wed, not Thu/Fri
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0

# Row 7 carries a custom row format; the new cells should use the default style
$ws.Range("E7:G7").Style = "Normal"

# Resize columns to fit their new contents (matches Excel's auto-fit result)
$ws.Columns(2).ColumnWidth = 8.330729166666666
$ws.Columns(3).ColumnWidth = 6.830729166666667
$ws.Columns(4).ColumnWidth = 3.9986979166666665
$ws.Columns(5).ColumnWidth = 4.166666666666667
$ws.Columns(6).ColumnWidth = 3.3307291666666665
$ws.Columns(7).ColumnWidth = 2.4986979166666665

$ws.Range("G7").Select()
